$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f9f10b82520>),
                ('model',
                 LogisticRegression(C=0.01, max_iter=1000, random_state=42,
                                    solver='saga'))])
'@

$ws.Range("B2").Value = 0.7333333333333333

$ws.Range("C2").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f9f10b82f70>, 'scaler': StandardScaler(), 'model__solver': 'saga', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 0.01}
'@

$ws.Range("D2").Value = 0.6602761322303702

$ws.Range("E2").Value = 0.635128205128205

$ws.Range("F2").Value = 0.8

$ws.Range("G2").Value = 0.6427771489809923

$ws.Range("H2").Value = 0.6317658730158729

$ws.Range("I2").Value = 0.6666666666666666

$ws.Range("J2").Value = 0.6975177304964539

$ws.Range("K2").Value = 0.6749999999999999

$ws.Range("L2").Value = 1

$ws.Range("N2").Value = @'
[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]
'@

$ws.Range("A3").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f9f901db1f0>),
                ('model',
                 LogisticRegression(C=0.01, max_iter=1000, random_state=42,
                                    solver='saga'))])
'@

$ws.Range("B3").Value = 0.732142857142857

$ws.Range("C3").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7fa79c082b20>, 'scaler': StandardScaler(), 'model__solver': 'saga', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 0.01}
'@

$ws.Range("D3").Value = 0.6920291399061833

$ws.Range("E3").Value = 0.6138057775557776

$ws.Range("G3").Value = 0.66629407892541

$ws.Range("H3").Value = 0.6736805555555555

$ws.Range("J3").Value = 0.7352836879432624

$ws.Range("K3").Value = 0.5958333333333332

$ws.Range("A4").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fa79c082cd0>),
                ('model',
                 LogisticRegression(C=0.01, max_iter=1000, random_state=42,
                                    solver='liblinear'))])
'@

$ws.Range("B4").Value = 0.7037373737373737

$ws.Range("C4").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f9f7034fe20>, 'scaler': StandardScaler(), 'model__solver': 'liblinear', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 0.01}
'@

$ws.Range("D4").Value = 0.6863842110669429

$ws.Range("E4").Value = 0.670060495060495

$ws.Range("F4").Value = 0.6206896551724138

$ws.Range("G4").Value = 0.6818331618905783

$ws.Range("H4").Value = 0.6828869047619048

$ws.Range("I4").Value = 0.9

$ws.Range("J4").Value = 0.702962962962963

$ws.Range("K4").Value = 0.6933333333333334

$ws.Range("L4").Value = 0.4736842105263158

$ws.Range("N4").Value = @'
[0 1 1 0 1 1 0 1 1 0 0 0 0 1 1 0 0 1 0 0 0 0 1 0]
'@

$ws.Range("A5").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f9f7034fd60>),
                ('model',
                 LogisticRegression(C=0.01, max_iter=1000, random_state=42,
                                    solver='saga'))])
'@

$ws.Range("B5").Value = 0.7499999999999999

$ws.Range("C5").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7fa6ccf7cca0>, 'scaler': StandardScaler(), 'model__solver': 'saga', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 0.01}
'@

$ws.Range("D5").Value = 0.6864114358775925

$ws.Range("E5").Value = 0.635421152921153

$ws.Range("G5").Value = 0.6711041734286211

$ws.Range("H5").Value = 0.6458796296296295

$ws.Range("J5").Value = 0.7304421768707482

$ws.Range("K5").Value = 0.6763888888888888

$ws.Range("A6").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fa6ccf7cdc0>),
                ('model',
                 LogisticRegression(C=0.01, max_iter=1000, random_state=42,
                                    solver='saga'))])
'@

$ws.Range("C6").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f9f10b85b80>, 'scaler': StandardScaler(), 'model__solver': 'saga', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 0.01}
'@

$ws.Range("D6").Value = 0.7549849134722101

$ws.Range("E6").Value = 0.7000421337921338

$ws.Range("G6").Value = 0.7186948995873639

$ws.Range("H6").Value = 0.6418088624338623

$ws.Range("J6").Value = 0.8121794871794871

$ws.Range("K6").Value = 0.7958333333333333
